# Regenerate the "K" column (col G) of the save_data sheet for lópez_reynaldo.
# The prior export mistakenly wrote a "Strike#" style count into column G;
# this regenerates the column using the corrected K (strikeouts) value per start.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-64 (row r corresponds to start index A(r-2)),
# in the same top-to-bottom order as the existing sheet rows.
$kValues = @(
    1,1,0,1,0,0,1,1,0,1,
    0,0,1,2,0,1,1,0,1,1,
    1,2,1,2,0,0,0,2,2,1,
    2,1,2,1,2,3,3,2,3,3,
    1,0,0,0,2,0,1,2,1,2,
    1,0,1,0,1,0,1,1,1,2,
    0,1,3
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
